$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Prueba"
$ws.Range("C1").Value = "Resultado"

$ws.Range("D1").Font.Underline = $true
$ws.Columns.Item(4).ColumnWidth = 23.140625

$ws.Range("D8").Select() | Out-Null
